# Fix Training Data Issue (#48)
# Data for 5-5-2007-08 was actually taken one day off due to the way
# NBA stats were shown -- correct the "Date" column (BF) from the
# "5-5-2007-08" placeholder text to the real ISO date "2008-05-05"
# for every data row, without letting Excel reinterpret the text as a
# serial date value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$dateCol = 58   # column BF

$rng = $ws.Range("BF$firstRow`:BF$lastRow")

# Force the range to a text format first so Excel does not silently
# convert the ISO-looking string "2008-05-05" into a date serial
# number when we assign it below.
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $dateCol).Value = "2008-05-05"
}

# Drop the temporary text formatting again so the cells end up with no
# explicit style applied, just like they started out.
$rng.ClearFormats()
